$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 280.07144
$ws.Range("I2").Value = 231.81818
$ws.Range("K2").Value = 231.81818
$ws.Range("M2").Value = -118.81818
$ws.Range("H17").Value = 11787.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 11787.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 35363.39999999999
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -35699.39999999999
$ws.Range("H42").Value = 150.25
$ws.Range("J42").Value = 185.6
$ws.Range("L42").Value = 556.8
$ws.Range("N42").Value = -1016.8
$ws.Range("H43").Value = 7634.85
$ws.Range("I43").Value = 3100.25
$ws.Range("K43").Value = 3100.25
$ws.Range("M43").Value = -3031.25
$ws.Range("H132").Value = 1735.8445
$ws.Range("I132").Value = 1634
$ws.Range("J132").Value = 2550.6
$ws.Range("K132").Value = 4902
$ws.Range("L132").Value = 7651.799999999999
$ws.Range("M132").Value = -2372
$ws.Range("N132").Value = -12711.8
$ws.Range("H137").Value = 4450.6
$ws.Range("I137").Value = 3677.875
$ws.Range("J137").Value = 4965.75
$ws.Range("K137").Value = 11033.625
$ws.Range("L137").Value = 14897.25
$ws.Range("M137").Value = -8483.625
$ws.Range("N137").Value = -19997.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1518.8
$ws.Range("J3").Value = 500
$ws.Range("L3").Value = 500
$ws.Range("N3").Value = -730
$ws.Range("H32").Value = 5014.8
$ws.Range("I32").Value = 3701.1462
$ws.Range("J32").Value = 10999.223
$ws.Range("K32").Value = 3701.1462
$ws.Range("L32").Value = 10999.223
$ws.Range("M32").Value = -3414.1462
$ws.Range("N32").Value = -11573.223
$ws.Range("H122").Value = 3524.1765
$ws.Range("I122").Value = 2638.875
$ws.Range("K122").Value = 7916.625
$ws.Range("M122").Value = -5466.625
$ws.Range("H124").Value = 59689.5
$ws.Range("J124").Value = 59689.5
$ws.Range("L124").Value = 59689.5
$ws.Range("N124").Value = -69509.5
$ws.Range("H125").Value = 75952.336
$ws.Range("J125").Value = 77542.8
$ws.Range("L125").Value = 77542.8
$ws.Range("N125").Value = -87382.8
$ws.Range("H132").Value = 2461.75
$ws.Range("I132").Value = 1742.0869
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 5226.2607
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -2696.2607
$ws.Range("N132").Value = -62102
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25415.58
$ws.Range("I31").Value = 3191.1714
$ws.Range("K31").Value = 3191.1714
$ws.Range("M31").Value = -2896.1714
$ws.Range("H34").Value = 25415.58
$ws.Range("I34").Value = 3191.1714
$ws.Range("K34").Value = 3191.1714
$ws.Range("M34").Value = -2989.1714
$ws.Range("H68").Value = 87994.60000000001
$ws.Range("J68").Value = 87994.60000000001
$ws.Range("L68").Value = 87994.60000000001
$ws.Range("N68").Value = -89492.60000000001
$ws.Range("H71").Value = 87994.60000000001
$ws.Range("J71").Value = 87994.60000000001
$ws.Range("L71").Value = 263983.8
$ws.Range("N71").Value = -271471.8
$ws.Range("H134").Value = 2908.8215
$ws.Range("I134").Value = 2170.5
$ws.Range("K134").Value = 6511.5
$ws.Range("M134").Value = -3976.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 52.142857
$ws.Range("I17").Value = 44
$ws.Range("J17").Value = 63
$ws.Range("K17").Value = 132
$ws.Range("L17").Value = 189
$ws.Range("M17").Value = 37
$ws.Range("N17").Value = -527
$ws.Range("H34").Value = 1774.8
$ws.Range("I34").Value = 1462.2
$ws.Range("J34").Value = 2400
$ws.Range("K34").Value = 4386.6
$ws.Range("L34").Value = 7200
$ws.Range("M34").Value = -4302.6
$ws.Range("N34").Value = -7368
$ws.Range("H39").Value = 1745.6
$ws.Range("J39").Value = 1000
$ws.Range("L39").Value = 3000
$ws.Range("N39").Value = -3588
$ws.Range("H55").Value = 1592.2142
$ws.Range("I55").Value = 1332.75
$ws.Range("J55").Value = 3149
$ws.Range("K55").Value = 3998.25
$ws.Range("L55").Value = 9447
$ws.Range("M55").Value = -3821.25
$ws.Range("N55").Value = -9801
$ws.Range("H88").Value = 19795.143
$ws.Range("I88").Value = 19775
$ws.Range("J88").Value = 19803.2
$ws.Range("K88").Value = 59325
$ws.Range("L88").Value = 59409.60000000001
$ws.Range("M88").Value = -58897
$ws.Range("N88").Value = -60265.60000000001
$ws.Range("H91").Value = 19795.143
$ws.Range("I91").Value = 19775
$ws.Range("J91").Value = 19803.2
$ws.Range("K91").Value = 59325
$ws.Range("L91").Value = 59409.60000000001
$ws.Range("M91").Value = -57843
$ws.Range("N91").Value = -62373.60000000001
$ws.Range("H141").Value = 10588.77
$ws.Range("I141").Value = 5647.25
$ws.Range("J141").Value = 12785
$ws.Range("K141").Value = 16941.75
$ws.Range("L141").Value = 38355
$ws.Range("M141").Value = -11761.75
$ws.Range("N141").Value = -48715
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1582.4
$ws.Range("I107").Value = 853.8333
$ws.Range("J107").Value = 2675.25
$ws.Range("K107").Value = 853.8333
$ws.Range("L107").Value = 2675.25
$ws.Range("M107").Value = 1066.1667
$ws.Range("N107").Value = -6515.25
$ws.Range("H111").Value = 46731
$ws.Range("J111").Value = 46731
$ws.Range("L111").Value = 46731
$ws.Range("N111").Value = -52865
$ws.Range("H113").Value = 2726.2917
$ws.Range("I113").Value = 1915.6666
$ws.Range("K113").Value = 1915.6666
$ws.Range("M113").Value = 254.3334
$ws.Range("H122").Value = 4904.512
$ws.Range("I122").Value = 4334.04
$ws.Range("J122").Value = 5795.875
$ws.Range("K122").Value = 13002.12
$ws.Range("L122").Value = 17387.625
$ws.Range("M122").Value = -10552.12
$ws.Range("N122").Value = -22287.625
$ws.Range("H132").Value = 6068.838
$ws.Range("I132").Value = 5368.6763
$ws.Range("J132").Value = 14004
$ws.Range("K132").Value = 16106.0289
$ws.Range("L132").Value = 42012
$ws.Range("M132").Value = -13576.0289
$ws.Range("N132").Value = -47072
$ws.Range("H134").Value = 54918
$ws.Range("J134").Value = 54918
$ws.Range("L134").Value = 164754
$ws.Range("N134").Value = -169824
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7030.857
$ws.Range("I16").Value = 3971.6667
$ws.Range("J16").Value = 12537.4
$ws.Range("K16").Value = 3971.6667
$ws.Range("L16").Value = 12537.4
$ws.Range("M16").Value = -3801.6667
$ws.Range("N16").Value = -12877.4
$ws.Range("H46").Value = 3542.3462
$ws.Range("I46").Value = 2287.625
$ws.Range("K46").Value = 2287.625
$ws.Range("M46").Value = -2099.625
$ws.Range("H55").Value = 2085380.8
$ws.Range("J55").Value = 4497.1113
$ws.Range("L55").Value = 4497.1113
$ws.Range("N55").Value = -4843.1113
$ws.Range("H100").Value = 7759.968
$ws.Range("I100").Value = 2460.6667
$ws.Range("J100").Value = 11106.895
$ws.Range("K100").Value = 2460.6667
$ws.Range("L100").Value = 11106.895
$ws.Range("M100").Value = -1919.6667
$ws.Range("N100").Value = -12188.895
$ws.Range("H132").Value = 8412.161
$ws.Range("I132").Value = 8133.706
$ws.Range("K132").Value = 24401.118
$ws.Range("M132").Value = -21871.118
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3503.838
$ws.Range("I132").Value = 2418.8965
$ws.Range("K132").Value = 7256.689499999999
$ws.Range("M132").Value = -4726.689499999999
$ws.Range("H136").Value = 3840.2856
$ws.Range("I136").Value = 3091.1667
$ws.Range("J136").Value = 8335
$ws.Range("K136").Value = 9273.500100000001
$ws.Range("L136").Value = 25005
$ws.Range("M136").Value = -6723.500100000001
$ws.Range("N136").Value = -30105
